$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new result row (row 57) produced by the latest script run.
# The leading apostrophe forces the date-looking string to stay text,
# matching the existing "Date" column cells (which are plain text, not
# real dates).
$ws.Range("A57").Value = "'2025-04-23"
$ws.Range("B57").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C57").Value = "NA"
$ws.Range("D57").Value = 1
